$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet: Overview
# Row2 now describes the ba774427 file, Row3 now describes the
# 90d90118 file (the two rows effectively swapped identity), and
# the 90d90118 row picks up a fresh "Ready for handoff" status.
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "ba774427-4f78-4031-ad1a-bc070f21edd8.md"
$wsOverview.Range("B2").Value = "e2e\ba774427-4f78-4031-ad1a-bc070f21edd8.md"

$wsOverview.Range("A3").Value = "90d90118-e242-44db-b2e6-e4a4b7d7e2d4.md"
$wsOverview.Range("B3").Value = "e2e\90d90118-e242-44db-b2e6-e4a4b7d7e2d4.md"
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-09-07 08:09:54"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/192bcff187947b0e07e5e386917f4000502f5300/e2e/90d90118-e242-44db-b2e6-e4a4b7d7e2d4.md", [Type]::Missing, [Type]::Missing, "e2e\ba774427-4f78-4031-ad1a-bc070f21edd8.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/192bcff187947b0e07e5e386917f4000502f5300/e2e/ba774427-4f78-4031-ad1a-bc070f21edd8.md", [Type]::Missing, [Type]::Missing, "e2e\90d90118-e242-44db-b2e6-e4a4b7d7e2d4.md") | Out-Null

# ---------------------------------------------------------------
# Sheet: zh-cn
# ---------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "ba774427-4f78-4031-ad1a-bc070f21edd8.md"
$wsZhCn.Range("G2").Value = "ba774427-4f78-4031-ad1a-bc070f21edd8.1b32bcab9903cb6e6644e907e97d1bfe70dbd639.zh-cn.xlf"
$wsZhCn.Range("I2").Value = "ba774427-4f78-4031-ad1a-bc070f21edd8.md"
$wsZhCn.Range("J2").Value = "ba774427-4f78-4031-ad1a-bc070f21edd8.1b32bcab9903cb6e6644e907e97d1bfe70dbd639.zh-cn.xlf"

$wsZhCn.Range("A3").Value = "90d90118-e242-44db-b2e6-e4a4b7d7e2d4.md"
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("G3").Value = "90d90118-e242-44db-b2e6-e4a4b7d7e2d4.246882f00340f1c95e596140032e420920c74481.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-09-07 08:09:42"
$wsZhCn.Range("I3").Value = "90d90118-e242-44db-b2e6-e4a4b7d7e2d4.md"
$wsZhCn.Range("J3").Value = "90d90118-e242-44db-b2e6-e4a4b7d7e2d4.246882f00340f1c95e596140032e420920c74481.zh-cn.xlf"
$wsZhCn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/192bcff187947b0e07e5e386917f4000502f5300/e2e/90d90118-e242-44db-b2e6-e4a4b7d7e2d4.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b2414d51438a2cb3e38ec85a5b8e4477dc628d5d/e2e/90d90118-e242-44db-b2e6-e4a4b7d7e2d4.md."

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/192bcff187947b0e07e5e386917f4000502f5300/e2e/90d90118-e242-44db-b2e6-e4a4b7d7e2d4.md", [Type]::Missing, [Type]::Missing, "ba774427-4f78-4031-ad1a-bc070f21edd8.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/314214a3afac75c4f127c34855b2340849a9bf7c/e2e/90d90118-e242-44db-b2e6-e4a4b7d7e2d4.md", [Type]::Missing, [Type]::Missing, "ba774427-4f78-4031-ad1a-bc070f21edd8.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/192bcff187947b0e07e5e386917f4000502f5300/e2e/ba774427-4f78-4031-ad1a-bc070f21edd8.md", [Type]::Missing, [Type]::Missing, "90d90118-e242-44db-b2e6-e4a4b7d7e2d4.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/314214a3afac75c4f127c34855b2340849a9bf7c/e2e/ba774427-4f78-4031-ad1a-bc070f21edd8.md", [Type]::Missing, [Type]::Missing, "90d90118-e242-44db-b2e6-e4a4b7d7e2d4.md") | Out-Null

$wsZhCn.Range("P1").ColumnWidth = 39.15

# ---------------------------------------------------------------
# Sheet: de-de
# ---------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "ba774427-4f78-4031-ad1a-bc070f21edd8.md"
$wsDeDe.Range("G2").Value = "ba774427-4f78-4031-ad1a-bc070f21edd8.1b32bcab9903cb6e6644e907e97d1bfe70dbd639.de-de.xlf"
$wsDeDe.Range("I2").Value = "ba774427-4f78-4031-ad1a-bc070f21edd8.md"
$wsDeDe.Range("J2").Value = "ba774427-4f78-4031-ad1a-bc070f21edd8.1b32bcab9903cb6e6644e907e97d1bfe70dbd639.de-de.xlf"

$wsDeDe.Range("A3").Value = "90d90118-e242-44db-b2e6-e4a4b7d7e2d4.md"
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("G3").Value = "90d90118-e242-44db-b2e6-e4a4b7d7e2d4.246882f00340f1c95e596140032e420920c74481.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-09-07 08:09:54"
$wsDeDe.Range("I3").Value = "90d90118-e242-44db-b2e6-e4a4b7d7e2d4.md"
$wsDeDe.Range("J3").Value = "90d90118-e242-44db-b2e6-e4a4b7d7e2d4.246882f00340f1c95e596140032e420920c74481.de-de.xlf"
$wsDeDe.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/192bcff187947b0e07e5e386917f4000502f5300/e2e/90d90118-e242-44db-b2e6-e4a4b7d7e2d4.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b2414d51438a2cb3e38ec85a5b8e4477dc628d5d/e2e/90d90118-e242-44db-b2e6-e4a4b7d7e2d4.md."

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/192bcff187947b0e07e5e386917f4000502f5300/e2e/90d90118-e242-44db-b2e6-e4a4b7d7e2d4.md", [Type]::Missing, [Type]::Missing, "ba774427-4f78-4031-ad1a-bc070f21edd8.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/83677caf5d7e6a5e691a9702d99023db283b9a9b/e2e/90d90118-e242-44db-b2e6-e4a4b7d7e2d4.md", [Type]::Missing, [Type]::Missing, "ba774427-4f78-4031-ad1a-bc070f21edd8.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/192bcff187947b0e07e5e386917f4000502f5300/e2e/ba774427-4f78-4031-ad1a-bc070f21edd8.md", [Type]::Missing, [Type]::Missing, "90d90118-e242-44db-b2e6-e4a4b7d7e2d4.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/83677caf5d7e6a5e691a9702d99023db283b9a9b/e2e/ba774427-4f78-4031-ad1a-bc070f21edd8.md", [Type]::Missing, [Type]::Missing, "90d90118-e242-44db-b2e6-e4a4b7d7e2d4.md") | Out-Null

$wsDeDe.Range("P1").ColumnWidth = 39.15
